$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL
$ws.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-status"

# Name
$ws.Range("B4").Value = "FRMedicationReconciliationStatus"

# Title
$ws.Range("B5").Value = "code system Interop'Santé - Statut d'une ligne de traitement d'une FCT"

# Date
$ws.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction
$ws.Range("B11").Value = "FRANCE"
